# Weekly data refresh: a new week's record is inserted right after the
# existing most-recent record (row 29), shifting all the older records
# down by one row. The new row is a duplicate of row 29's content, and
# row 29 then gets the brand-new date (one week after the previous most
# recent date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (blank) row at position 30; rows 30..59 shift to 31..60.
$ws.Rows.Item(30).Insert()

# Populate the freshly inserted row 30 with the same content row 29 had
# (row 29 is the template for the new weekly record).
$ws.Cells.Item(30, 1).Value2  = 9
$ws.Cells.Item(30, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(30, 3).Value2  = "Metropolitana"
$ws.Cells.Item(30, 4).Value2  = 44623
$ws.Cells.Item(30, 5).Value2  = 13
$ws.Cells.Item(30, 6).Value2  = 100112029
$ws.Cells.Item(30, 7).Value2  = "Orégano"
$ws.Cells.Item(30, 8).Value2  = "Sin especificar"
$ws.Cells.Item(30, 9).Value2  = "Primera"
$ws.Cells.Item(30, 10).Value2 = 16
$ws.Cells.Item(30, 11).Value2 = 16000
$ws.Cells.Item(30, 12).Value2 = 16000
$ws.Cells.Item(30, 13).Value2 = 16000
$ws.Cells.Item(30, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(30, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(30, 16).Value2 = 5333
$ws.Cells.Item(30, 17).Value2 = 3
$ws.Cells.Item(30, 18).Value2 = "Hortaliza"

# Match the date cell's number format (date/time) to the rest of column D.
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat

# The previously-most-recent row (29) now carries the new week's date.
$ws.Cells.Item(29, 4).Value2 = 44664
